$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.121.49"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "4.026.00"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.700"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.86%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.749"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000325"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "4.672.80"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "4.037.28"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D20").Value = "72.032.18"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "98.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.73%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +26.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "681.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.427"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "0.0₃0824"
$ws.Range("E39").Value = "  -9.06%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.32%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0486"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.151"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.48%  "
$ws.Range("E46").Value = "  +4.75%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.86%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
